$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row for "Shalom Dube" / "Business level 3 first year " (row 8).
# All rows below it shift up by one; the table/autofilter range auto-adjusts.
$ws.Rows.Item(8).Delete()

# Normalize the free-text "What year are you in" values (column C) onto a
# consistent small set of labels.
$ws.Cells.Item(5, 3).Value2 = "College"
$ws.Cells.Item(8, 3).Value2 = "Final year of college"
$ws.Cells.Item(18, 3).Value2 = "2nd year of college"
$ws.Cells.Item(21, 3).Value2 = "2nd year of college"
$ws.Cells.Item(28, 3).Value2 = "1st year of college"
$ws.Cells.Item(29, 3).Value2 = "1st year of college"
$ws.Cells.Item(31, 3).Value2 = "2nd year of college "

# Keep the current selection on the last data row, matching the observed
# post-edit cursor position.
$ws.Range("C34").Select()
